$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date field
$ws.Range("G4").Value = [DateTime]::FromOADate(44076.78246170176)

# Customer name
$ws.Range("G7").Value = "FMU"

# Comments
$ws.Range("F10").Value = "Let's see what is it"

# Row 19 - Digital (GB)
$ws.Range("F19").Value = 3000
$ws.Range("G19").Value = 14
$ws.Range("H19").Value = 40320

# Row 22 - Online Storage (GB)
$ws.Range("F22").Value = 3000
$ws.Range("H22").Value = 1152

# Row 24 - Registration fee
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 200
$ws.Range("H24").Value = 200

# Row 25 - AWA contribution / Entity
$ws.Range("E25").Value = "public"
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 500
$ws.Range("H25").Value = 500

# Row 26 - Management fee (per year)
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 60
$ws.Range("H26").Value = 60

# Row 27 - Storage (reels/per year) / Period (years)
# "Period (years)" is entered as text "25" (a text label rather than a
# numeric quantity), while F27 (Qty column) is a true number. Use a scratch
# cell to build the text value, then paste just the value into E27 so the
# cell keeps its original (General) style/format instead of picking up a
# Text number-format.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "25"
$scratch.Copy()
$ws.Range("E27").PasteSpecial(-4163) # xlPasteValues
$scratch.Clear()

$ws.Range("F27").Value = 25
$ws.Range("G27").Value = 40
$ws.Range("H27").Value = 25000

# Row 29 - piqlReader
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 79900
$ws.Range("H29").Value = 79900

# Row 30 - Installation and training
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 3000
$ws.Range("H30").Value = 3000

# Row 31 - Service agreement (per year) / Type
$ws.Range("E31").Value = "gold"
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 2500
$ws.Range("H31").Value = 2500

# Row 32 - Shipment cost / Reels
$ws.Range("E32").Value = 25
$ws.Range("G32").Value = 20
$ws.Range("H32").Value = 500

# Row 33 - TOTAL
$ws.Range("H33").Value = 162952

# Row 34 - Total to pay from the second term
$ws.Range("H34").Value = 12532
